$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.305.31"
$ws.Range("E2").Value = "  +1.75%  "

$ws.Range("D3").Value = "2.309.49"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'302.23"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").Value = "'100.92"
$ws.Range("E6").Value = "  +5.52%  "

$ws.Range("D7").Value = "'0.503"
$ws.Range("E7").Value = "  +0.68%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  +5.34%  "

$ws.Range("D10").Value = "'36.70"
$ws.Range("E10").Value = "  +10.37%  "

$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("D12").Value = "'18.59"
$ws.Range("E12").Value = "  +12.46%  "

$ws.Range("E13").Value = "  +1.62%  "

$ws.Range("D14").Value = "'6.97"
$ws.Range("E14").Value = "  +3.61%  "

$ws.Range("D15").Value = "2.671.58"
$ws.Range("E15").Value = "  +1.11%  "

$ws.Range("D16").Value = "2.335.17"
$ws.Range("E16").Value = "  +1.64%  "

$ws.Range("D17").Value = "'0.805"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "43.168.79"
$ws.Range("E18").Value = "  +1.66%  "

$ws.Range("D19").Value = "'12.74"
$ws.Range("E19").Value = "  +11.53%  "

$ws.Range("D20").Value = "'6.23"
$ws.Range("E20").Value = "  +4.46%  "

$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +1.96%  "

$ws.Range("D22").Value = "'68.24"

$ws.Range("D23").Value = "'237.05"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("D24").Value = "'2.25"
$ws.Range("E24").Value = "  +15.05%  "

$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").Value = "'25.19"
$ws.Range("E27").Value = "  +4.12%  "

$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'35.07"
$ws.Range("E28").Value = "  +4.27%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'168.50"
$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.07"
$ws.Range("E30").Value = "  -4.07%  "

$ws.Range("D31").Value = "'9.21"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").Value = "'5.04"
$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("D34").Value = "'17.86"
$ws.Range("E34").Value = "  +5.41%  "

$ws.Range("D35").Value = "'4.68"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").Value = "'2.42"
$ws.Range("E36").Value = "  +0.38%  "

$ws.Range("D37").Value = "'0.0699"
$ws.Range("E37").Value = "  +1.41%  "

$ws.Range("D38").Value = "'2.84"
$ws.Range("E38").Value = "  +1.81%  "

$ws.Range("D39").Value = "'1.80"
$ws.Range("E39").Value = "  +4.00%  "

$ws.Range("D40").Value = "'0.101"
$ws.Range("E40").Value = "  +1.58%  "

$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").Value = "'2.30"
$ws.Range("E42").Value = "  -2.82%  "

$ws.Range("D43").Value = "1.991.81"
$ws.Range("E43").Value = "  +2.00%  "

$ws.Range("D44").Value = "'0.0291"
$ws.Range("E44").Value = "  +4.48%  "

$ws.Range("D45").Value = "'10.17"
$ws.Range("E45").Value = "  +5.00%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.92"
$ws.Range("E46").Value = "  +4.36%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'17.72"
$ws.Range("E47").Value = "  +1.68%  "

$ws.Range("D48").Value = "'55.91"
$ws.Range("E48").Value = "  +6.67%  "

$ws.Range("D49").Value = "'1.56"
$ws.Range("E49").Value = "  +5.93%  "

$ws.Range("D50").Value = "2.538.75"
$ws.Range("E50").Value = "  +1.08%  "

$ws.Range("D51").Value = "'4.54"
$ws.Range("E51").Value = "  +0.40%  "
